# Insert a new weekly price record at row 88 ("Hortaliza, Vega Monumental
# Concepción - Alcachofa"). Inserting the row shifts the existing rows
# 88-91 down to 89-92, which reproduces the rest of the diff automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 88, pushing rows 88..91 down to 89..92.
$ws.Rows.Item(88).Insert()

# Populate the new row 88 with the new weekly record.
$ws.Range("A88").Value = 11
$ws.Range("B88").Value = "Vega Monumental Concepción"
$ws.Range("C88").Value = "Bíobío"
$ws.Range("D88").Value = 45106
$ws.Range("E88").Value = 8
$ws.Range("F88").Value = 100112013
$ws.Range("G88").Value = "Alcachofa"
$ws.Range("H88").Value = "Argentina(o)"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 100
$ws.Range("K88").Value = 14000
$ws.Range("L88").Value = 15000
$ws.Range("M88").Value = 14500
$ws.Range("N88").Value = "`$/caja 50 unidades"
$ws.Range("O88").Value = "Provincia de Limarí"
$ws.Range("P88").Value = 290
$ws.Range("Q88").Value = 50
$ws.Range("R88").Value = "Hortaliza"
